$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biweekly_Gas_Samples")
$ws.Range("A1").Value = "test"
